# Update "想去人数" (want-to-go count) figures in F column across sheets,
# matching the regenerated data snapshot from the gh-pages build.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 2552
$wsExhibit.Range("F12").Value = 7691
$wsExhibit.Range("F20").Value = 9449
$wsExhibit.Range("F37").Value = 1491
$wsExhibit.Range("F40").Value = 219

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 55
$wsShow.Range("F22").Value = 35

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 55
$wsAll.Range("F7").Value = 2552
$wsAll.Range("F18").Value = 7691
$wsAll.Range("F24").Value = 9449
$wsAll.Range("F37").Value = 1491
$wsAll.Range("F41").Value = 219
